$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Order products alphabetically and " -> "Order products alphabetically "
#    (drop the trailing "and"; the leftover single space stays)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Order products alphabetically and ", $true, $false, $false, $false, $false, $true, 1, $false, "Order products alphabetically ", 2)

# Locate that paragraph and split its single run into: "O" | "rder products alphabetically" | " "
$pStart = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Order products alphabetically ")) {
        $pStart = $p.Range.Start
    }
}
if ($pStart -ne $null) {
    # boundary after "O" (1 char)
    $rSplit1 = $d.Range($pStart, $pStart + 1)
    $rSplit1.Bold = 1
    $rSplit1.Bold = 0
    # boundary before the trailing " " (last char, position 29..30)
    $rSplit2 = $d.Range($pStart + 29, $pStart + 30)
    $rSplit2.Bold = 1
    $rSplit2.Bold = 0
}

# ---------------------------------------------------------------------------
# 2) "Order products reverse alphabetically" paragraph is untouched.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 3) "Order products in increasing and decreasing order of price"
#    -> "Order products in order of price"
#    and move the hidden "_GoBack" bookmark from the end of the
#    "Sort and Filter products" paragraph to between "Order products in "
#    and "order of price".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Order products in increasing and decreasing order of price", $true, $false, $false, $false, $false, $true, 1, $false, "Order products in order of price", 2)

$pStart2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Order products in order of price")) {
        $pStart2 = $p.Range.Start
    }
}
if ($pStart2 -ne $null) {
    # split into "Order products in " (18 chars) | "order of price" (14 chars)
    $rSplit3 = $d.Range($pStart2 + 18, $pStart2 + 18 + 14)
    $rSplit3.Bold = 1
    $rSplit3.Bold = 0
}

# Remove the old bookmark location and recreate it at the new split point.
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

$bmPos = $pStart2 + 18
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 4) "Filter by a single product tag at a time" -> "Filter by a single product tag "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Filter by a single product tag at a time", $true, $false, $false, $false, $false, $true, 1, $false, "Filter by a single product tag ", 2)
